$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Phase 0: give the new date cells (B3:C19) the same date style as B2:C2 ---
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Phase 1: prime shared-string table in exact target order (41..76) ---
$ws.Cells.Item(3, "D").Value = "San Diego, CA"
$ws.Cells.Item(3, "E").Value = "American Society of Safety Engineers"
$ws.Cells.Item(3, "A").Value = "Conducting a Safety Audit"
$ws.Cells.Item(3, "H").Value = "File"
$ws.Cells.Item(4, "E").Value = "Georgia Public Safety Training Center"
$ws.Cells.Item(5, "A").Value = "Basic Course for Health Assessment and Consultation"
$ws.Cells.Item(5, "E").Value = "Agency for Toxic Substances and Disease Registry "
$ws.Cells.Item(6, "A").Value = "Introduction to Risk Assessment Guidance for Superfund"
$ws.Cells.Item(6, "E").Value = "US Environmental Protection Agency"
$ws.Cells.Item(7, "A").Value = "Time Series and Forecasting"
$ws.Cells.Item(7, "E").Value = "Practical Stats, Inc."
$ws.Cells.Item(8, "A").Value = "Untangling Multivariate Relationships"
$ws.Cells.Item(9, "A").Value = "Nondetects and Data Analysis, Statistical Methods for Censored Environmental Data"
$ws.Cells.Item(10, "A").Value = "Principles of Quality Assurance and Quality Control in Environmental Field Programs"
$ws.Cells.Item(10, "D").Value = "Montgomery, AL"
$ws.Cells.Item(10, "E").Value = "Northwest Environmental Training Center"
$ws.Cells.Item(11, "A").Value = "Exposure Risk Training"
$ws.Cells.Item(11, "E").Value = "US Environmental Protection Agency Office of Air Quality Planning and Standards"
$ws.Cells.Item(12, "A").Value = "Intermediate Incident Command System for Expanding Incidents (ICS300)"
$ws.Cells.Item(12, "E").Value = "US Centers for Disease Control and Prevention"
$ws.Cells.Item(13, "A").Value = "Applied Environmental Statistics"
$ws.Cells.Item(4, "A").Value = "Advanced Command System for Command and General Staff, Complex Incidents and MACS (ICS400)"
$ws.Cells.Item(14, "A").Value = "Sampling for Hazardous Matrials"
$ws.Cells.Item(15, "A").Value = "Air Monitoring for Emergency Response"
$ws.Cells.Item(16, "A").Value = "Interpretation of Biomonitoring Data using Physiologically Based Pharmokinetic (PBPK) Modeling"
$ws.Cells.Item(16, "D").Value = "RTP, NC"
$ws.Cells.Item(17, "A").Value = "Environmental Sampling, Sampling Reliability and Data Quality Objectives for the Health Assessment Process"
$ws.Cells.Item(17, "E").Value = "Envirostat, Inc."
$ws.Cells.Item(16, "E").Value = "CIIT Centers for Health Research"
$ws.Cells.Item(18, "A").Value = "AERMOD Air Dispersion Modeling"
$ws.Cells.Item(18, "D").Value = "Chicago, IL"
$ws.Cells.Item(18, "E").Value = "Lakes Environmental"
$ws.Cells.Item(19, "D").Value = "Edison, NJ"
$ws.Cells.Item(19, "A").Value = "Health and Safety (40 Hour HAZWOPER)"
$ws.Cells.Item(14, "D").Value = "Longmont, CO"
$ws.Cells.Item(4, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"

# --- Phase 2: fill remaining cells (duplicated strings + all numeric values) ---
$ws.Cells.Item(3, "B").Value = 36617
$ws.Cells.Item(3, "C").Value = 36617
$ws.Cells.Item(4, "B").Value = 40241
$ws.Cells.Item(4, "C").Value = 405483
$ws.Cells.Item(4, "D").Value = "Atlanta, GA"
$ws.Cells.Item(5, "B").Value = 40805
$ws.Cells.Item(5, "C").Value = 40809
$ws.Cells.Item(5, "D").Value = "Atlanta, GA"
$ws.Cells.Item(5, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(6, "B").Value = 41345
$ws.Cells.Item(6, "C").Value = 41347
$ws.Cells.Item(6, "D").Value = "Atlanta, GA"
$ws.Cells.Item(6, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(7, "B").Value = 41345
$ws.Cells.Item(7, "C").Value = 41346
$ws.Cells.Item(7, "D").Value = "Atlanta, GA"
$ws.Cells.Item(7, "F").Value = 13
$ws.Cells.Item(7, "G").Value = "Hours"
$ws.Cells.Item(7, "H").Value = "File"
$ws.Cells.Item(8, "B").Value = 40982
$ws.Cells.Item(8, "C").Value = 40983
$ws.Cells.Item(8, "D").Value = "Atlanta, GA"
$ws.Cells.Item(8, "E").Value = "Practical Stats, Inc."
$ws.Cells.Item(8, "F").Value = 14
$ws.Cells.Item(8, "G").Value = "Hours"
$ws.Cells.Item(8, "H").Value = "File"
$ws.Cells.Item(9, "B").Value = 40771
$ws.Cells.Item(9, "C").Value = 40772
$ws.Cells.Item(9, "D").Value = "Atlanta, GA"
$ws.Cells.Item(9, "E").Value = "Practical Stats, Inc."
$ws.Cells.Item(9, "F").Value = 14
$ws.Cells.Item(9, "G").Value = "Hours"
$ws.Cells.Item(9, "H").Value = "File"
$ws.Cells.Item(10, "B").Value = 40277
$ws.Cells.Item(10, "C").Value = 40278
$ws.Cells.Item(10, "F").Value = 13
$ws.Cells.Item(10, "G").Value = "Hours"
$ws.Cells.Item(10, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(11, "B").Value = 42234
$ws.Cells.Item(11, "C").Value = 42236
$ws.Cells.Item(11, "D").Value = "Atlanta, GA"
$ws.Cells.Item(11, "F").Value = 19
$ws.Cells.Item(11, "G").Value = "Hours"
$ws.Cells.Item(12, "B").Value = 40118
$ws.Cells.Item(12, "C").Value = 40120
$ws.Cells.Item(12, "D").Value = "Atlanta, GA"
$ws.Cells.Item(12, "H").Value = "File"
$ws.Cells.Item(13, "B").Value = 40287
$ws.Cells.Item(13, "C").Value = 40291
$ws.Cells.Item(13, "D").Value = "Atlanta, GA"
$ws.Cells.Item(13, "E").Value = "Practical Stats, Inc."
$ws.Cells.Item(14, "B").Value = 39322
$ws.Cells.Item(14, "C").Value = 39324
$ws.Cells.Item(14, "E").Value = "US Environmental Protection Agency"
$ws.Cells.Item(14, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(15, "B").Value = 38880
$ws.Cells.Item(15, "C").Value = 38881
$ws.Cells.Item(15, "D").Value = "Atlanta, GA"
$ws.Cells.Item(15, "E").Value = "US Environmental Protection Agency"
$ws.Cells.Item(15, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(16, "B").Value = 38985
$ws.Cells.Item(16, "C").Value = 38989
$ws.Cells.Item(16, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(17, "B").Value = 39629
$ws.Cells.Item(17, "C").Value = 39632
$ws.Cells.Item(17, "D").Value = "Atlanta, GA"
$ws.Cells.Item(17, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(18, "B").Value = 39352
$ws.Cells.Item(18, "C").Value = 39353
$ws.Cells.Item(18, "H").Value = "\\cdc.gov\private\M309\Hzd3\training\Certificates"
$ws.Cells.Item(19, "B").Value = 39818
$ws.Cells.Item(19, "C").Value = 39822
$ws.Cells.Item(19, "E").Value = "US Environmental Protection Agency"

# --- Phase 3: selection / active sheet bookkeeping ---
$ws.Activate() | Out-Null
$ws.Range("A20").Select() | Out-Null
